# "more budget estimation, routing sheet"
# Fill in direct-cost values for the "Non-recurring" budget section (rows 6-8),
# fill in zeroes for the insurance/security line items (rows 12-13),
# extend the H28 "Recurring" total formula to include F12 and F13,
# and add a new "first time direct cost" total row (row 31) that sums
# the Non-recurring (H26) and Recurring (H28) totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New direct-cost figures under the Non-recurring section
$ws.Range("F6").Value = 500
$ws.Range("F7").Value = 30
$ws.Range("F8").Value = 500

# Insurance / security costs, currently both nil
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 0

# Recurring total now also needs to add in the insurance/security cells
$ws.Range("H28").Formula = "=F5+F11+F12+F13+F19"

# New summary row: "first time direct cost" = Non-recurring total + Recurring total
$ws.Range("A31").Value = "first time direct cost"
$ws.Range("H31").Formula = "=H26+H28"

# Keep the active selection in sync with where Excel would have left the
# cursor after adding the new row (one row below the new content).
$ws.Range("H32").Select()
